$d = $word.ActiveDocument

# 1. Add new bullet "Mario level small khi nhấn xuống thì rớt khỏi màn hình"
#    right after the "Chuyển level của mario thì bị rớt khỏi màn hình." bullet
#    (still inside the BUG section, before "Những thứ cần làm:").
$rng1 = $d.Content
$rng1.Find.Execute("Chuyển level của mario thì bị rớt khỏi màn hình.")
$rng1.Collapse(0)
$rng1.InsertAfter("`rMario level small khi nhấn xuống thì rớt khỏi màn hình")

# 2. Turn the trailing empty bullet into a new "CLASS ColorBrick" heading,
#    followed by a bullet describing the class.
$p = $d.Paragraphs.Last
$rng2 = $p.Range
$rng2.InsertAfter("CLASS ColorBrick`r Là class để mấy khung màu sắc tại vị trí phía trên cùng")
$p.Style = "Heading1"
